$d = $word.ActiveDocument

# 1-4: Replace "Bytem: [[ADDRESS_n]].[[PHONE_n]]" with "Bytem: [[ADDRESS_n]].: +420 [[AMOUNT_n]]"
for ($i = 1; $i -le 4; $i++) {
    $old = "Bytem: [[ADDRESS_$i]].[[PHONE_$i]]"
    $new = "Bytem: [[ADDRESS_$i]].: +420 [[AMOUNT_$i]]"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# 5: Replace "Měsíční nájemné činí 18 000 Kč." with "Měsíční nájemné činí [[AMOUNT_5]]."
$d.Content.Find.Execute("Měsíční nájemné činí 18 000 Kč.", $true, $false, $false, $false, $false, $true, 1, $false, "Měsíční nájemné činí [[AMOUNT_5]].", 2)

# 6: Replace "Kauce činí 36 000 Kč a bude uhrazena nejpozději při podpisu této smlouvy." with "Kauce činí [[AMOUNT_6]] a bude uhrazena nejpozději při podpisu této smlouvy."
$d.Content.Find.Execute("Kauce činí 36 000 Kč a bude uhrazena nejpozději při podpisu této smlouvy.", $true, $false, $false, $false, $false, $true, 1, $false, "Kauce činí [[AMOUNT_6]] a bude uhrazena nejpozději při podpisu této smlouvy.", 2)
